# Auto-generated edit script: updates crypto price/volume data
# per commit "Updated cryptos list on Wed Sep 27 10:56:22 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold plain text in the source data (e.g. "213.91", "0.0610").
# Some of the new values look like ordinary numbers to Excel; left alone, Excel would
# silently convert them to floating point numbers (dropping significant trailing zeros
# such as "227.00" -> 227, or introducing binary floating-point noise). Force those
# specific cells to Text format first so the literal string is preserved exactly.
$textFormatCells = @("D5", "D9", "D10", "D11", "D16", "D18", "D20", "D22", "D23", "D25", "D29", "D34", "D37", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "26.406.90"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.614.50"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "213.91"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.0610"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "1.842.09"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "1.636.28"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "64.64"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "26.422.97"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "227.00"
$ws.Range("E18").Value = "  +6.40%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "9.10"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "15.37"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "1.448.89"
$ws.Range("E33").Value = "  +8.73%  "
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "0.558"
$ws.Range("E37").Value = "  -4.94%  "
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "1.753.61"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("D44").Value = "0.763"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "61.96"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "0.911"
$ws.Range("E46").Value = "  -10.79%  "
$ws.Range("D47").Value = "87.85"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.49"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0960"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  +1.97%  "
